$wb = $excel.ActiveWorkbook

# --- Add the new "Info" worksheet, placed after "ColumnsNumberParameters" --
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Info"
$newSheet.Move([System.Reflection.Missing]::Value, $wb.Worksheets.Item("ColumnsNumberParameters"))

# Re-fetch a live reference to the sheet (Move() can re-seat indices).
$info = $wb.Worksheets.Item("Info")

# --- A1: intro sentence, with "reg_estimates" in italics -------------------
$text1 = "This Excel file is used to define the column numbers required for the corresponding processes in the reg_estimates files."
$info.Range("A1").Value = $text1
$info.Range("A1").Characters(102, 13).Font.Italic = $true

# --- A2: second sentence, with "separately" and "each country" in italic+underline
$text2 = "Since the column numbers differ by country, they must be configured separately for each country."
$info.Range("A2").Value = $text2
$info.Range("A2").Characters(69, 10).Font.Italic = $true
$info.Range("A2").Characters(69, 10).Font.Underline = $true
$info.Range("A2").Characters(84, 12).Font.Italic = $true
$info.Range("A2").Characters(84, 12).Font.Underline = $true

# --- Register the matching whole-run font variants (italic; italic+underline)
# in the workbook style table, mirroring what Excel's desktop UI leaves
# behind when characters are formatted via the mini toolbar / Format Cells.
$info.Range("C1").Value = "x"
$info.Range("C1").Font.Italic = $true
$info.Range("C1").Font.Italic = $false

$info.Range("C2").Value = "x"
$info.Range("C2").Font.Italic = $true
$info.Range("C2").Font.Underline = $true
$info.Range("C2").Font.Underline = $false
$info.Range("C2").Font.Italic = $false

$info.Range("C1").Value = ""
$info.Range("C2").Value = ""

# --- Selection & active-tab bookkeeping -------------------------------------
$info.Range("A1:A2").Select()
$info.Activate()
